$wb = $excel.ActiveWorkbook

# --- Sheet "CÔ DIỄM" (index 1 / sheet1.xml) ---
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 24, shifting everything (incl. formulas) down by one.
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 (H/I/K columns only)
$ws.Cells.Item(24, 8).Value = "30/06/2024"
$ws.Cells.Item(24, 9).Value = 10
$ws.Cells.Item(24, 11).Value = "chưa làm giấy"

# Append two new rows at the bottom (now rows 125 and 126)
$ws.Cells.Item(125, 1).Value = "23/06/2024"
$ws.Cells.Item(125, 2).Value = "Duy lấy tiền lời ngày 30,31 = 12tr"
$ws.Cells.Item(125, 3).Value = 12000
$ws.Cells.Item(125, 4).Formula = "=D124+C125"

$ws.Cells.Item(126, 2).Value = "Duy cho cô Diễm vay 10tr"
$ws.Cells.Item(126, 3).Value = -10000
$ws.Cells.Item(126, 4).Formula = "=D125+C126"

# Make "CÔ DIỄM" the active/selected sheet and set the view/selection.
$ws.Activate()
$ws.Range("K24").Select()

# --- Sheet "Tổng Kết" (index 5 / sheet5.xml) - no longer the selected tab ---
$tk = $wb.Worksheets.Item(5)
$tk.Range("B21").Select()

# Re-activate "CÔ DIỄM" so it ends up as the active sheet/tab.
$ws.Activate()
